# feat: add support for loading card weight from data
#
# Adds a new "Weight" column (O) to the card data sheet:
#   - O1 header = "Weight"
#   - O2 = 1, O3 = 1, O4 = 2 (numeric weight per card row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, same header style as the rest of row 1 (col N)
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Weight"

# New numeric "Weight" values for the three data rows
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("O4").Value = 2

# Leave the user's selection on the newly-entered cell, matching the final
# edit position after typing the new column's data, and scroll the view so
# the new column is visible.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O4").Select()
